$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "product_id" header in B1 with "quantity", and clear the old
# quantity header cell in C1 so the product_id column is effectively removed.
$ws.Range("B1").Value = "quantity"
$ws.Range("C1").Clear()

# Update the active selection to match the target state
$ws.Range("F8").Select()
